$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Email Money Transfers (AP)")

# --- Remove the sample contact rows (names + emails + mailto hyperlinks) ---
$ws.Hyperlinks.Delete()
$ws.Range("B3:B5").Clear()
$ws.Range("C3:C8").Clear()

# D3:D9 keep their blank "contact name" formatting but lose the sample values
$ws.Range("D3:D9").ClearContents()
$ws.Range("D3").Copy()
$ws.Range("D5:D9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Replace the security question / answer sample text ---
$ws.Range("G3:G16").Value = "What is the three-letter acronym of the Graduate Student Society?"
$ws.Range("H3:I16").Value = "GSS"

# --- Email validation now starts at row 3 (sample rows gone) ---
$emailRange = $ws.Range("C3:C1048576")
$emailRange.Validation.Delete()
$emailRange.Validation.Add(7, 1, 5, "=SEARCH(""."",C3,(SEARCH(""@"",C3,1))+2)>0")
$emailRange.Validation.InputTitle = "Contact email"
$emailRange.Validation.InputMessage = "Enter valid email"

# --- Column G needs to be wide enough for the longer question text ---
$ws.Columns.Item(7).ColumnWidth = 77.6640625

# --- Update the active selection to the new security-question column ---
$ws.Activate()
$ws.Range("G3:G16").Select()
